$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Plain text-safe updates -------------------------------------------------
# These values are not ambiguous with Excel's automatic number/date
# detection (they either contain extra separators like "58.332.89", are
# percentage strings with surrounding spaces, or are plainly non-numeric),
# so a direct Range.Value assignment keeps them stored as text, exactly as
# they were before the edit.
$textUpdates = @{
    "D2" = '58.332.89'
    "E2" = '  -0.52%  '
    "D3" = '3.141.10'
    "E3" = '  +1.21%  '
    "E4" = '  +0.01%  '
    "E5" = '  +1.15%  '
    "E6" = '  -0.97%  '
    "E7" = '  +0.06%  '
    "D8" = '3.141.75'
    "E8" = '  +1.24%  '
    "E9" = '  +0.68%  '
    "E10" = '  -2.71%  '
    "E11" = '  +0.25%  '
    "E12" = '  +1.83%  '
    "D13" = '3.681.58'
    "E13" = '  +1.45%  '
    "E14" = '  +3.37%  '
    "E15" = '  -4.79%  '
    "E16" = '  -0.67%  '
    "D17" = '58.349.23'
    "E17" = '  -0.52%  '
    "D18" = '3.141.54'
    "E18" = '  +1.43%  '
    "E19" = '  -0.85%  '
    "E20" = '  -1.01%  '
    "E21" = '  -1.79%  '
    "E22" = '  +0.28%  '
    "E23" = '  -0.02%  '
    "E24" = '  +1.29%  '
    "E25" = '  +2.61%  '
    "E26" = '  -0.74%  '
    "E27" = '  -0.14%  '
    "E28" = '  +1.59%  '
    "B29" = 'InternetComputer(DFINITY)'
    "C29" = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
    "E29" = '  +2.00%  '
    "B30" = 'USDe'
    "C30" = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
    "E30" = '  +0.03%  '
    "E31" = '  -3.58%  '
    "E32" = '  +1.31%  '
    "E33" = '  +0.20%  '
    "E34" = '  -1.38%  '
    "E35" = '  +3.09%  '
    "E36" = '  +2.26%  '
    "E37" = '  +2.32%  '
    "E38" = '  -2.48%  '
    "E39" = '  -5.27%  '
    "E40" = '  +11.86%  '
    "E41" = '  -1.57%  '
    "E42" = '  +5.65%  '
    "E43" = '  +2.30%  '
    "D44" = '3.181.82'
    "E44" = '  +1.33%  '
    "E45" = '  -0.41%  '
    "E46" = '  +0.05%  '
    "E47" = '  +2.59%  '
    "D48" = '2.277.33'
    "E48" = '  -0.30%  '
    "E49" = '  +4.54%  '
    "E50" = '  -1.80%  '
    "E51" = '  +1.50%  '
}

foreach ($ref in $textUpdates.Keys) {
    $ws.Range($ref).Value = $textUpdates[$ref]
}

# --- Numeric-looking text updates --------------------------------------------
# These "Price" values (e.g. "532.81", "0.392", "7.40") look like plain
# numbers, so a direct Range.Value assignment would make Excel silently
# convert them to real numbers (losing the trailing zeros / exact text
# seen in the source data, e.g. "7.40" -> 7.4). To keep them as literal
# text - matching the original inline-string cells - each one is entered
# as a quoted string-literal formula and then pasted back as a value,
# which preserves Excel's text storage without altering cell formatting.
$numericTextUpdates = @{
    "D5" = '532.81'
    "D6" = '142.38'
    "D10" = '7.14'
    "D12" = '0.392'
    "D15" = '25.64'
    "D19" = '6.11'
    "D20" = '12.80'
    "D22" = '343.25'
    "D25" = '67.64'
    "D29" = '7.40'
    "D30" = '0.999'
    "D33" = '21.10'
    "D35" = '4.80'
    "D36" = '157.88'
    "D37" = '6.22'
    "D38" = '26.35'
    "D45" = '36.70'
    "D50" = '20.65'
}

foreach ($ref in $numericTextUpdates.Keys) {
    $val = $numericTextUpdates[$ref]
    $cell = $ws.Range($ref)
    $cell.Formula = '="' + $val + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

$excel.CutCopyMode = 0
